# IATI partner activities - troubleshoot missing iati.cloud activities
# Adds a new "Notes" column and a new row documenting an activity that
# is not being returned as a linked activity in the iati.cloud extraction.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IATI activity IDs")

# --- New row 21: a J-PAL / Acumen MECS-linked activity that's missing ---
# Fill left-to-right so new shared strings land in the same order the
# original author typed them (A, C, E), reusing existing strings for
# B (extending org) and D (fund) exactly as used in the other FCDO rows.
$ws.Range("A21").Value = "US-EIN-042103594-GCCI-3978870"
$ws.Range("B21").Value = "Foreign, Commonwealth and Development Office"
$ws.Range("C21").Value = "GB-GOV-1-300049"
$ws.Range("D21").Value = "FCDO Research - Programmes"
$ws.Range("E21").Value = "J-PAL"

# Highlight the iati_id in dark red to flag it for follow-up.
$ws.Range("A21").Font.Color = 192

# --- New "Notes" column, documenting the troubleshooting finding ---
$ws.Range("F1").Value = "Notes"
# Match the same header formatting (bold white text, navy fill, centred)
# used by the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F21").Value = "not being returned as a linked activity in iati.cloud extraction"

$ws.Range("C12").Select()
